# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 4124.9375
$ws.Range("I51").Value = 3369.8462
$ws.Range("M51").Value = -2885.8462
$ws.Range("K51").Value = 3369.8462
# Row 52
$ws.Range("J52").Value = 10000
$ws.Range("I52").Value = 664.5
$ws.Range("H52").Value = 6265.8
$ws.Range("M52").Value = -1833.5
$ws.Range("N52").Value = -30320
$ws.Range("K52").Value = 1993.5
$ws.Range("L52").Value = 30000
# Row 62
$ws.Range("M62").Value = -5768.6924
$ws.Range("N62").Value = -9348
$ws.Range("K62").Value = 6392.6924
$ws.Range("J62").Value = 8100
$ws.Range("L62").Value = 8100
$ws.Range("I62").Value = 6392.6924
$ws.Range("H62").Value = 6712.8125
# Row 65
$ws.Range("K65").Value = 31963.462
$ws.Range("M65").Value = -28843.462
$ws.Range("I65").Value = 6392.6924
$ws.Range("H65").Value = 6712.8125
$ws.Range("J65").Value = 8100
$ws.Range("N65").Value = -46740
$ws.Range("L65").Value = 40500
# Row 99
$ws.Range("I99").Value = 255.75
$ws.Range("J99").Value = 845
$ws.Range("H99").Value = 452.16666
$ws.Range("N99").Value = -5531
$ws.Range("M99").Value = 730.75
$ws.Range("K99").Value = 767.25
$ws.Range("L99").Value = 2535
# Row 125
$ws.Range("K125").Value = 13752
$ws.Range("H125").Value = 7094755.5
$ws.Range("M125").Value = -11292
$ws.Range("I125").Value = 1528
# Row 131
$ws.Range("J131").Value = 7993.3335
$ws.Range("H131").Value = 6016.1
$ws.Range("M131").Value = -10466.1432
$ws.Range("L131").Value = 23980.0005
$ws.Range("K131").Value = 15506.1432
$ws.Range("I131").Value = 5168.7144
$ws.Range("N131").Value = -34060.00049999999

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3087128.5
$ws.Range("M2").Value = -3472697
$ws.Range("I2").Value = 3472810
$ws.Range("K2").Value = 3472810
# Row 61
$ws.Range("K61").Value = 1695.125
$ws.Range("M61").Value = -1483.125
$ws.Range("I61").Value = 1695.125
$ws.Range("H61").Value = 2632.7896
# Row 102
$ws.Range("I102").Value = 6947773.5
$ws.Range("K102").Value = 6947773.5
$ws.Range("H102").Value = 5212186
$ws.Range("M102").Value = -6946151.5
# Row 116
$ws.Range("M116").Value = -3470516
$ws.Range("I116").Value = 3472810
$ws.Range("K116").Value = 3472810
$ws.Range("H116").Value = 3087128.5
# Row 136
$ws.Range("I136").Value = 1695.125
$ws.Range("H136").Value = 2632.7896
$ws.Range("M136").Value = -2535.375
$ws.Range("K136").Value = 5085.375

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("I3").Value = 3472810
$ws.Range("H3").Value = 3087128.5
$ws.Range("K3").Value = 3472810
$ws.Range("M3").Value = -3472696

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("M16").Value = -1752.7778
$ws.Range("K16").Value = 2039.7778
$ws.Range("H16").Value = 2321
$ws.Range("I16").Value = 2039.7778
# Row 58
$ws.Range("I58").Value = 3647.9546
$ws.Range("N58").Value = -3564.8333
$ws.Range("H58").Value = 3543.1428
$ws.Range("K58").Value = 3647.9546
$ws.Range("J58").Value = 3158.8333
$ws.Range("M58").Value = -3444.9546
$ws.Range("L58").Value = 3158.8333
# Row 105
$ws.Range("I105").Value = 2831.5833
$ws.Range("K105").Value = 2831.5833
$ws.Range("J105").Value = 3037.75
$ws.Range("H105").Value = 2883.125
$ws.Range("L105").Value = 3037.75
$ws.Range("M105").Value = -1084.5833
$ws.Range("N105").Value = -6531.75
# Row 113
$ws.Range("M113").Value = 130.2221999999999
$ws.Range("K113").Value = 2039.7778
$ws.Range("H113").Value = 2321
$ws.Range("I113").Value = 2039.7778
# Row 122
$ws.Range("L122").Value = 11595
$ws.Range("K122").Value = 10031.4999
$ws.Range("M122").Value = -7581.499899999999
$ws.Range("I122").Value = 3343.8333
$ws.Range("J122").Value = 3865
$ws.Range("H122").Value = 3517.5557
$ws.Range("N122").Value = -16495
# Row 132
$ws.Range("I132").Value = 60365.65
$ws.Range("K132").Value = 181096.95
$ws.Range("J132").Value = 82965.91
$ws.Range("M132").Value = -178566.95
$ws.Range("N132").Value = -253957.73
$ws.Range("L132").Value = 248897.73
$ws.Range("H132").Value = 69244.32000000001
# Row 136
$ws.Range("I136").Value = 3647.9546
$ws.Range("L136").Value = 9476.499899999999
$ws.Range("H136").Value = 3543.1428
$ws.Range("K136").Value = 10943.8638
$ws.Range("N136").Value = -14576.4999
$ws.Range("J136").Value = 3158.8333
$ws.Range("M136").Value = -8393.863799999999

$ws = $wb.Worksheets.Item("CUL")
# Row 120
$ws.Range("H120").Value = 12082.909
$ws.Range("J120").Value = 20546.6
$ws.Range("N120").Value = -71315.79999999999
$ws.Range("L120").Value = 61639.8

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("I102").Value = 6174226.5
$ws.Range("K102").Value = 6174226.5
$ws.Range("H102").Value = 5229906.5
$ws.Range("M102").Value = -6172604.5
# Row 122
$ws.Range("L122").Value = 6813.8181
$ws.Range("K122").Value = 1540323.48
$ws.Range("M122").Value = -1537873.48
$ws.Range("I122").Value = 513441.16
$ws.Range("J122").Value = 2271.2727
$ws.Range("H122").Value = 337726.5
$ws.Range("N122").Value = -11713.8181

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("K7").Value = 3500
$ws.Range("I7").Value = 3500
$ws.Range("M7").Value = -3388
$ws.Range("H7").Value = 6624.75
# Row 16
$ws.Range("M16").Value = -458.9286
$ws.Range("K16").Value = 628.9286
$ws.Range("H16").Value = 855.9143
$ws.Range("I16").Value = 628.9286
# Row 40
$ws.Range("N40").Value = -11605.667
$ws.Range("K40").Value = 6596.7144
$ws.Range("J40").Value = 11333.667
$ws.Range("H40").Value = 8017.8
$ws.Range("M40").Value = -6460.7144
$ws.Range("L40").Value = 11333.667
$ws.Range("I40").Value = 6596.7144
# Row 41
$ws.Range("I41").Value = 18500
$ws.Range("J41").Value = 44999.5
$ws.Range("K41").Value = 18500
$ws.Range("N41").Value = -45875.5
$ws.Range("H41").Value = 31749.75
$ws.Range("M41").Value = -18062
$ws.Range("L41").Value = 44999.5
# Row 126
$ws.Range("K126").Value = 10500
$ws.Range("H126").Value = 6624.75
$ws.Range("M126").Value = -8030
$ws.Range("I126").Value = 3500
# Row 132
$ws.Range("I132").Value = 3084.186
$ws.Range("K132").Value = 9252.558000000001
$ws.Range("H132").Value = 3736.0386
$ws.Range("M132").Value = -6722.558000000001
# Row 136
$ws.Range("I136").Value = 77413.516
$ws.Range("H136").Value = 63031.293
$ws.Range("M136").Value = -229690.548
$ws.Range("K136").Value = 232240.548

$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 34949.5
$ws.Range("J80").Value = 34949.5
$ws.Range("N80").Value = -36945.5
$ws.Range("L80").Value = 34949.5
# Row 81
$ws.Range("I81").Value = 23811582
$ws.Range("J81").Value = 1653.3334
$ws.Range("K81").Value = 47623164
$ws.Range("H81").Value = 16668603
$ws.Range("M81").Value = -47622103
$ws.Range("N81").Value = -5428.6668
$ws.Range("L81").Value = 3306.6668
# Row 83
$ws.Range("J83").Value = 34949.5
$ws.Range("L83").Value = 104848.5
$ws.Range("N83").Value = -114832.5
$ws.Range("H83").Value = 34949.5
# Row 84
$ws.Range("J84").Value = 1653.3334
$ws.Range("L84").Value = 16533.334
$ws.Range("M84").Value = -238110516
$ws.Range("H84").Value = 16668603
$ws.Range("I84").Value = 23811582
$ws.Range("K84").Value = 238115820
$ws.Range("N84").Value = -27141.334
# Row 88
$ws.Range("H88").Value = 24000
$ws.Range("J88").Value = 24000
$ws.Range("L88").Value = 24000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("N88").Value = -24812
$ws.Range("M88").ClearContents()
# Row 91
$ws.Range("L91").Value = 24000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 24000
$ws.Range("K91").Value = 0
$ws.Range("H91").Value = 24000
$ws.Range("N91").Value = -26808
$ws.Range("M91").ClearContents()
# Row 126
$ws.Range("L126").Value = 5173.125
$ws.Range("J126").Value = 1724.375
$ws.Range("H126").Value = 2589.3333
$ws.Range("N126").Value = -10113.125
